$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of facility row ranges (in column AA, "average_county_temperature")
# to the new NOAA-derived temperature values, per the commit:
# "Updated temperature with NOAA data" / "Added back NAICS 311230" / "Added merged datasets"

$ws.Range("AA2:AA7").Value = 13.75752314814816
$ws.Range("AA20:AA31").Value = 13.75752314814816
$ws.Range("AA32:AA37").Value = 3.38888888888889
$ws.Range("AA44:AA60").Value = 12.93898809523811
$ws.Range("AA73:AA90").Value = 19.79629629629628
$ws.Range("AA91:AA96").Value = 13.75752314814816
$ws.Range("AA97:AA114").Value = 3.38888888888889
